$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price and Volume(1h) data columns keep their text (string) type,
# so values such as "1.00" or "25.13" are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '59.919.88'
$ws.Range("E2").Value = '  -6.16%  '
$ws.Range("D3").Value = '2.988.62'
$ws.Range("E3").Value = '  -6.46%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '570.43'
$ws.Range("E5").Value = '  -4.15%  '
$ws.Range("D6").Value = '125.21'
$ws.Range("E6").Value = '  -8.79%  '
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("D8").Value = '2.984.91'
$ws.Range("E8").Value = '  -6.45%  '
$ws.Range("E9").Value = '  -2.61%  '
$ws.Range("E10").Value = '  -9.45%  '
$ws.Range("D11").Value = '5.05'
$ws.Range("E11").Value = '  -6.05%  '
$ws.Range("E12").Value = '  -4.13%  '
$ws.Range("E13").Value = '  -9.60%  '
$ws.Range("E14").Value = '  -7.22%  '
$ws.Range("E15").Value = '  +0.83%  '
$ws.Range("D16").Value = '3.481.67'
$ws.Range("E16").Value = '  -6.33%  '
$ws.Range("D17").Value = '2.986.02'
$ws.Range("E17").Value = '  -6.32%  '
$ws.Range("D18").Value = '59.948.87'
$ws.Range("E18").Value = '  -6.08%  '
$ws.Range("D19").Value = '6.44'
$ws.Range("E19").Value = '  -2.32%  '
$ws.Range("D20").Value = '424.98'
$ws.Range("E20").Value = '  -8.44%  '
$ws.Range("D21").Value = '13.08'
$ws.Range("E21").Value = '  -6.57%  '
$ws.Range("E22").Value = '  -4.99%  '
$ws.Range("E23").Value = '  -8.61%  '
$ws.Range("D24").Value = '12.93'
$ws.Range("E24").Value = '  -2.59%  '
$ws.Range("E25").Value = '  -4.98%  '
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("E28").Value = '  -6.28%  '
$ws.Range("D29").Value = '1.94'
$ws.Range("E29").Value = '  -7.71%  '
$ws.Range("E30").Value = '  -7.44%  '
$ws.Range("D31").Value = '6.16'
$ws.Range("E31").Value = '  -10.53%  '
$ws.Range("D32").Value = '25.13'
$ws.Range("D33").Value = '0.0937'
$ws.Range("E33").Value = '  -8.35%  '
$ws.Range("E34").Value = '  -5.29%  '
$ws.Range("D35").Value = '0.926'
$ws.Range("E35").Value = '  -9.98%  '
$ws.Range("D36").Value = '49.99'
$ws.Range("E36").Value = '  -3.43%  '
$ws.Range("D37").Value = '2.05'
$ws.Range("E37").Value = '  -17.17%  '
$ws.Range("E38").Value = '  -11.73%  '
$ws.Range("D39").Value = '8.35'
$ws.Range("E39").Value = '  +1.76%  '
$ws.Range("E40").Value = '  -10.45%  '
$ws.Range("E41").Value = '  -5.42%  '
$ws.Range("D42").Value = '375.61'
$ws.Range("E42").Value = '  -6.03%  '
$ws.Range("D43").Value = '2.665.90'
$ws.Range("E43").Value = '  -5.07%  '
$ws.Range("E44").Value = '  -8.66%  '
$ws.Range("E46").Value = '  -7.94%  '
$ws.Range("D47").Value = '120.14'
$ws.Range("E47").Value = '  -6.97%  '
$ws.Range("D48").Value = '1.99'
$ws.Range("E48").Value = '  -7.64%  '
$ws.Range("E49").Value = '  -4.39%  '
$ws.Range("D50").Value = '23.58'
$ws.Range("E50").Value = '  -8.26%  '
$ws.Range("E51").Value = '  -8.13%  '
